$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 83729.5
$ws.Range("I12").Value = 400
$ws.Range("J12").Value = 167059
$ws.Range("K12").Value = 400
$ws.Range("L12").Value = 167059
$ws.Range("M12").Value = -230
$ws.Range("N12").Value = -167399
$ws.Range("H19").Value = 1207.9474
$ws.Range("I19").Value = 849.6667
$ws.Range("J19").Value = 1822.1428
$ws.Range("K19").Value = 849.6667
$ws.Range("L19").Value = 1822.1428
$ws.Range("M19").Value = -674.6667
$ws.Range("N19").Value = -2172.1428
$ws.Range("H69").Value = 4636
$ws.Range("J69").Value = 4480
$ws.Range("L69").Value = 13440
$ws.Range("N69").Value = -15188
$ws.Range("H72").Value = 4636
$ws.Range("J72").Value = 4480
$ws.Range("L72").Value = 40320
$ws.Range("N72").Value = -49056
$ws.Range("H74").Value = 3779.9
$ws.Range("I74").Value = 3374.75
$ws.Range("J74").Value = 4050
$ws.Range("K74").Value = 3374.75
$ws.Range("L74").Value = 4050
$ws.Range("M74").Value = -2438.75
$ws.Range("N74").Value = -5922
$ws.Range("H77").Value = 3779.9
$ws.Range("I77").Value = 3374.75
$ws.Range("J77").Value = 4050
$ws.Range("K77").Value = 16873.75
$ws.Range("L77").Value = 20250
$ws.Range("M77").Value = -12193.75
$ws.Range("N77").Value = -29610
$ws.Range("H98").Value = 678.34784
$ws.Range("I98").Value = 645.5454999999999
$ws.Range("K98").Value = 645.5454999999999
$ws.Range("M98").Value = 852.4545000000001
$ws.Range("H122").Value = 678.34784
$ws.Range("I122").Value = 645.5454999999999
$ws.Range("K122").Value = 1936.6365
$ws.Range("M122").Value = 513.3635000000002
$ws.Range("H129").Value = 3695.4856
$ws.Range("J129").Value = 921.0714
$ws.Range("L129").Value = 2763.2142
$ws.Range("N129").Value = -12763.2142
$ws.Range("H137").Value = 2043.8
$ws.Range("I137").Value = 2262.5715
$ws.Range("J137").Value = 1852.375
$ws.Range("K137").Value = 6787.7145
$ws.Range("L137").Value = 5557.125
$ws.Range("M137").Value = -4237.7145
$ws.Range("N137").Value = -10657.125
$ws.Range("H138").Value = 3460.0986
$ws.Range("I138").Value = 2222.5454
$ws.Range("J138").Value = 4015.7346
$ws.Range("K138").Value = 6667.6362
$ws.Range("L138").Value = 12047.2038
$ws.Range("M138").Value = -1527.6362
$ws.Range("N138").Value = -22327.2038

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32335.172
$ws.Range("I32").Value = 14326.571
$ws.Range("J32").Value = 171001.4
$ws.Range("K32").Value = 14326.571
$ws.Range("L32").Value = 171001.4
$ws.Range("M32").Value = -14039.571
$ws.Range("N32").Value = -171575.4
$ws.Range("H74").Value = 1579.3334
$ws.Range("I74").Value = 1506.6666
$ws.Range("J74").Value = 1797.3334
$ws.Range("K74").Value = 1506.6666
$ws.Range("L74").Value = 1797.3334
$ws.Range("M74").Value = -632.6666
$ws.Range("N74").Value = -3545.3334
$ws.Range("H77").Value = 1579.3334
$ws.Range("I77").Value = 1506.6666
$ws.Range("J77").Value = 1797.3334
$ws.Range("K77").Value = 7533.333000000001
$ws.Range("L77").Value = 8986.666999999999
$ws.Range("M77").Value = -3165.333000000001
$ws.Range("N77").Value = -17722.667
$ws.Range("H110").Value = 50106388
$ws.Range("I110").Value = 62632656
$ws.Range("K110").Value = 62632656
$ws.Range("M110").Value = -62630611
$ws.Range("H122").Value = 2307.8823
$ws.Range("I122").Value = 2011
$ws.Range("K122").Value = 6033
$ws.Range("M122").Value = -3583
$ws.Range("H132").Value = 10820.766
$ws.Range("I132").Value = 12389.537
$ws.Range("J132").Value = 2349.4
$ws.Range("K132").Value = 37168.611
$ws.Range("L132").Value = 7048.200000000001
$ws.Range("M132").Value = -34638.611
$ws.Range("N132").Value = -12108.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 8602.799999999999
$ws.Range("I25").Value = 2003.5
$ws.Range("J25").Value = 35000
$ws.Range("K25").Value = 2003.5
$ws.Range("L25").Value = 35000
$ws.Range("M25").Value = -1768.5
$ws.Range("N25").Value = -35470
$ws.Range("H120").Value = 33853.5
$ws.Range("J120").Value = 33853.5
$ws.Range("L120").Value = 33853.5
$ws.Range("N120").Value = -43529.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29892.416
$ws.Range("I31").Value = 1507.7778
$ws.Range("J31").Value = 59368.77
$ws.Range("K31").Value = 1507.7778
$ws.Range("L31").Value = 59368.77
$ws.Range("M31").Value = -1212.7778
$ws.Range("N31").Value = -59958.77
$ws.Range("H34").Value = 29892.416
$ws.Range("I34").Value = 1507.7778
$ws.Range("J34").Value = 59368.77
$ws.Range("K34").Value = 1507.7778
$ws.Range("L34").Value = 59368.77
$ws.Range("M34").Value = -1305.7778
$ws.Range("N34").Value = -59772.77
$ws.Range("H80").Value = 13087
$ws.Range("J80").Value = 13087
$ws.Range("L80").Value = 13087
$ws.Range("N80").Value = -15333
$ws.Range("H83").Value = 13087
$ws.Range("J83").Value = 13087
$ws.Range("L83").Value = 39261
$ws.Range("N83").Value = -50493
$ws.Range("H122").Value = 1178.2
$ws.Range("J122").Value = 1197.75
$ws.Range("L122").Value = 3593.25
$ws.Range("N122").Value = -8493.25
$ws.Range("H132").Value = 2960.7144
$ws.Range("I132").Value = 2858.8
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 8576.400000000001
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -6046.400000000001
$ws.Range("N132").Value = -20057
$ws.Range("H134").Value = 1094.3055
$ws.Range("I134").Value = 659.9231
$ws.Range("J134").Value = 2223.7
$ws.Range("K134").Value = 1979.7693
$ws.Range("L134").Value = 6671.099999999999
$ws.Range("M134").Value = 555.2307000000001
$ws.Range("N134").Value = -11741.1
$ws.Range("H135").Value = 49117.8
$ws.Range("J135").Value = 49117.8
$ws.Range("L135").Value = 49117.8
$ws.Range("N135").Value = -59257.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4504.909
$ws.Range("I56").Value = 4504.909
$ws.Range("K56").Value = 4504.909
$ws.Range("M56").Value = -3974.909
$ws.Range("H123").Value = 3143.3333
$ws.Range("J123").Value = 4250
$ws.Range("L123").Value = 12750
$ws.Range("N123").Value = -17650
$ws.Range("H131").Value = 618476.4399999999
$ws.Range("I131").Value = 590
$ws.Range("J131").Value = 704294
$ws.Range("K131").Value = 1770
$ws.Range("L131").Value = 2112882
$ws.Range("M131").Value = 3270
$ws.Range("N131").Value = -2122962
$ws.Range("H136").Value = 3052.8572
$ws.Range("I136").Value = 2597.5
$ws.Range("J136").Value = 3660
$ws.Range("K136").Value = 7792.5
$ws.Range("L136").Value = 10980
$ws.Range("M136").Value = -2692.5
$ws.Range("N136").Value = -21180
$ws.Range("H137").Value = 40449.934
$ws.Range("I137").Value = 103403
$ws.Range("J137").Value = 8973.4
$ws.Range("K137").Value = 310209
$ws.Range("L137").Value = 26920.2
$ws.Range("M137").Value = -305109
$ws.Range("N137").Value = -37120.2
$ws.Range("H138").Value = 14328.75
$ws.Range("I138").Value = 18438.334
$ws.Range("J138").Value = 2000
$ws.Range("K138").Value = 55315.00199999999
$ws.Range("L138").Value = 6000
$ws.Range("M138").Value = -50175.00199999999
$ws.Range("N138").Value = -16280

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2635.9167
$ws.Range("I102").Value = 2013.5
$ws.Range("J102").Value = 3258.3333
$ws.Range("K102").Value = 2013.5
$ws.Range("L102").Value = 3258.3333
$ws.Range("M102").Value = -391.5
$ws.Range("N102").Value = -6502.3333
$ws.Range("H126").Value = 4588.125
$ws.Range("I126").Value = 4448.5
$ws.Range("K126").Value = 13345.5
$ws.Range("M126").Value = -10875.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1362.3143
$ws.Range("I55").Value = 1841.875
$ws.Range("K55").Value = 1841.875
$ws.Range("M55").Value = -1668.875
$ws.Range("H68").Value = 1672.3334
$ws.Range("I68").Value = 1135
$ws.Range("J68").Value = 1941
$ws.Range("K68").Value = 1135
$ws.Range("L68").Value = 1941
$ws.Range("M68").Value = -386
$ws.Range("N68").Value = -3439
$ws.Range("H71").Value = 1672.3334
$ws.Range("I71").Value = 1135
$ws.Range("J71").Value = 1941
$ws.Range("K71").Value = 5675
$ws.Range("L71").Value = 9705
$ws.Range("M71").Value = -1931
$ws.Range("N71").Value = -17193
$ws.Range("H122").Value = 4301.278
$ws.Range("I122").Value = 3565.3333
$ws.Range("J122").Value = 5773.1665
$ws.Range("K122").Value = 10695.9999
$ws.Range("L122").Value = 17319.4995
$ws.Range("M122").Value = -8245.999899999999
$ws.Range("N122").Value = -22219.4995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 25647358
$ws.Range("I62").Value = 38466536
$ws.Range("J62").Value = 9000
$ws.Range("K62").Value = 38466536
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = -38465912
$ws.Range("N62").Value = -10248
$ws.Range("H65").Value = 25647358
$ws.Range("I65").Value = 38466536
$ws.Range("J65").Value = 9000
$ws.Range("K65").Value = 192332680
$ws.Range("L65").Value = 45000
$ws.Range("M65").Value = -192329560
$ws.Range("N65").Value = -51240
$ws.Range("H81").Value = 200548.1
$ws.Range("J81").Value = 333863
$ws.Range("L81").Value = 667726
$ws.Range("N81").Value = -669848
$ws.Range("H84").Value = 200548.1
$ws.Range("J84").Value = 333863
$ws.Range("L84").Value = 3338630
$ws.Range("N84").Value = -3349238
$ws.Range("H96").Value = 76924500
$ws.Range("I96").Value = 111112670
$ws.Range("J96").Value = 1102
$ws.Range("K96").Value = 111112670
$ws.Range("L96").Value = 1102
$ws.Range("M96").Value = -111111297
$ws.Range("N96").Value = -3848
$ws.Range("H122").Value = 2027.4667
$ws.Range("I122").Value = 2284.3333
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 6852.999899999999
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -4402.999899999999
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 21184.283
$ws.Range("I132").Value = 2415.4736
$ws.Range("J132").Value = 68731.92999999999
$ws.Range("K132").Value = 7246.4208
$ws.Range("L132").Value = 206195.79
$ws.Range("M132").Value = -4716.4208
$ws.Range("N132").Value = -211255.79
